# "Generate Report for Handoff"
# Updates the localization-status report: file "b.md" now has a new handoff
# (status "Ready for handoff") with a freshly generated xlf handoff file and
# handoff datetime, for both the zh-cn and de-de locales, plus the summary
# status shown on the Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the "b.md" file; update its rolled-up status ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 is "b.md" ---
$zhcnHandoffAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92c9f916633a84ed8c7d2ec184640d4c2402d284/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-08 12:22:46"

$zhcn.Hyperlinks.Add($zhcn.Range("C3"), $zhcnHandoffAddr, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null

# --- de-de sheet: row 3 is "b.md" ---
$dedeHandoffAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20d8efa39a471f3499661dbd7762100e95daa345/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-08 12:22:51"

$dede.Hyperlinks.Add($dede.Range("C3"), $dedeHandoffAddr, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
